$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing correction factor values for E13 and E21
$ws.Range("E13").Value = 0.3160978559
$ws.Range("E21").Value = 0.14254470129999999

# Update the active selection to E21 to match the saved view state
$ws.Range("E21").Select()
